$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.16289592760181
$ws.Range("C2").Value = 0.5972850678733032
$ws.Range("J2").Value = 0.03167420814479638
$ws.Range("P2").Value = 0.1266968325791855
$ws.Range("S2").Value = 0.08144796380090498
$ws.Range("C3").Value = 0.03597122302158273
$ws.Range("J3").Value = 0.05035971223021583
$ws.Range("P3").Value = 0.7194244604316546
$ws.Range("S3").Value = 0.1942446043165468
$ws.Range("J4").Value = 0.06521739130434782
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.04273504273504274
$ws.Range("D6").Value = 0.008547008547008548
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.04700854700854701
$ws.Range("J6").Value = 0.188034188034188
$ws.Range("O6").Value = 0.0170940170940171
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.4273504273504273
$ws.Range("B7").Value = 0.09947643979057591
$ws.Range("D7").Value = 0.02094240837696335
$ws.Range("E7").Value = 0.01047120418848168
$ws.Range("F7").Value = 0.07853403141361257
$ws.Range("J7").Value = 0.09424083769633508
$ws.Range("O7").Value = 0.02094240837696335
$ws.Range("Q7").Value = 0.1675392670157068
$ws.Range("R7").Value = 0.08900523560209424
$ws.Range("S7").Value = 0.418848167539267
$ws.Range("B8").Value = 0.09151785714285714
$ws.Range("D8").Value = 0.02008928571428572
$ws.Range("E8").Value = 0.002232142857142857
$ws.Range("F8").Value = 0.06919642857142858
$ws.Range("J8").Value = 0.1026785714285714
$ws.Range("O8").Value = 0.02455357142857143
$ws.Range("Q8").Value = 0.1629464285714286
$ws.Range("R8").Value = 0.1205357142857143
$ws.Range("S8").Value = 0.40625
$ws.Range("B9").Value = 0.08176100628930817
$ws.Range("D9").Value = 0.02515723270440252
$ws.Range("F9").Value = 0.06918238993710692
$ws.Range("J9").Value = 0.1132075471698113
$ws.Range("O9").Value = 0.01257861635220126
$ws.Range("Q9").Value = 0.1949685534591195
$ws.Range("R9").Value = 0.1006289308176101
$ws.Range("S9").Value = 0.4025157232704403
$ws.Range("B10").Value = 0.08561341571050309
$ws.Range("D10").Value = 0.02383053839364519
$ws.Range("F10").Value = 0.08737864077669903
$ws.Range("J10").Value = 0.1032656663724625
$ws.Range("O10").Value = 0.02383053839364519
$ws.Range("Q10").Value = 0.1800529567519859
$ws.Range("R10").Value = 0.116504854368932
$ws.Range("S10").Value = 0.3795233892321271
$ws.Range("G11").Value = 0.1114864864864865
$ws.Range("J11").Value = 0.09797297297297297
$ws.Range("K11").Value = 0.1824324324324324
$ws.Range("L11").Value = 0.5878378378378378
$ws.Range("S11").Value = 0.02027027027027027
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1988636363636364
$ws.Range("K12").Value = 0.005681818181818182
$ws.Range("L12").Value = 0.01136363636363636
$ws.Range("S12").Value = 0.03409090909090909
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.0125
$ws.Range("H15").Value = 0.1541666666666667
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.3375
$ws.Range("K15").Value = 0.05833333333333333
$ws.Range("M15").Value = 0.02083333333333333
$ws.Range("N15").Value = 0.004166666666666667
$ws.Range("O15").Value = 0.09166666666666666
$ws.Range("S15").Value = 0.2583333333333334
$ws.Range("F16").Value = 0.02515723270440252
$ws.Range("H16").Value = 0.1509433962264151
$ws.Range("I16").Value = 0.08176100628930817
$ws.Range("J16").Value = 0.3836477987421384
$ws.Range("K16").Value = 0.1006289308176101
$ws.Range("M16").Value = 0.02515723270440252
$ws.Range("N16").Value = 0.01257861635220126
$ws.Range("O16").Value = 0.0440251572327044
$ws.Range("S16").Value = 0.1761006289308176
$ws.Range("F17").Value = 0.02680965147453083
$ws.Range("H17").Value = 0.2010723860589812
$ws.Range("I17").Value = 0.07774798927613941
$ws.Range("J17").Value = 0.3994638069705094
$ws.Range("K17").Value = 0.09919571045576407
$ws.Range("M17").Value = 0.008042895442359249
$ws.Range("O17").Value = 0.04557640750670242
$ws.Range("S17").Value = 0.1420911528150134
$ws.Range("F18").Value = 0.02439024390243903
$ws.Range("H18").Value = 0.1829268292682927
$ws.Range("I18").Value = 0.0975609756097561
$ws.Range("J18").Value = 0.4024390243902439
$ws.Range("K18").Value = 0.08536585365853659
$ws.Range("M18").Value = 0.02845528455284553
$ws.Range("O18").Value = 0.05691056910569105
$ws.Range("S18").Value = 0.1219512195121951
$ws.Range("F19").Value = 0.01482701812191104
$ws.Range("H19").Value = 0.2232289950576606
$ws.Range("I19").Value = 0.06589785831960461
$ws.Range("J19").Value = 0.3500823723228995
$ws.Range("K19").Value = 0.1243822075782537
$ws.Range("M19").Value = 0.0214168039538715
$ws.Range("N19").Value = 0.0008237232289950577
$ws.Range("O19").Value = 0.08484349258649095
$ws.Range("S19").Value = 0.114497528830313
